$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246 (pushes existing rows 246..300 down to 247..301)
$ws.Rows("246:246").Insert()

# Fill in the new row 246 with the new weekly record
$ws.Range("A246").Value = 3
$ws.Range("B246").Value = "Femacal de La Calera"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = Get-Date -Year 2023 -Month 9 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("E246").Value = 5
$ws.Range("F246").Value = 100112026
$ws.Range("G246").Value = "Haba"
$ws.Range("H246").Value = "Sin especificar"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 40
$ws.Range("K246").Value = 15000
$ws.Range("L246").Value = 15000
$ws.Range("M246").Value = 15000
$ws.Range("N246").Value = "$/saco 25 kilos"
$ws.Range("O246").Value = "Provincia de Limarí"
$ws.Range("P246").Value = 600
$ws.Range("Q246").Value = 25
$ws.Range("R246").Value = "Hortaliza"
